$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.817481054725187
$ws.Range("D2").Value = 8.956687002177262
$ws.Range("E2").Value = 13.59163149916536
$ws.Range("F2").Value = 33.84103078573688
$ws.Range("G2").Value = 3.660266244116293
$ws.Range("J2").Value = 10.12650920219991
$ws.Range("K2").Value = 11.17907715981412
$ws.Range("M2").Value = 15.71016208632966
$ws.Range("N2").Value = 19.66633514129256
$ws.Range("O2").Value = 25.57303292257555
$ws.Range("B3").Value = 7.74649179613065
$ws.Range("D3").Value = 8.939537045287773
$ws.Range("E3").Value = 13.59980362885076
$ws.Range("F3").Value = 33.87928389297875
$ws.Range("G3").Value = 3.662244156708274
$ws.Range("J3").Value = 10.14894479144957
$ws.Range("K3").Value = 10.90073912465009
$ws.Range("M3").Value = 15.61018128412094
$ws.Range("N3").Value = 19.72757654451025
$ws.Range("O3").Value = 25.63100130601634
$ws.Range("B4").Value = 7.704416437718603
$ws.Range("D4").Value = 8.930330081366453
$ws.Range("E4").Value = 13.60704204713831
$ws.Range("F4").Value = 33.91050717482582
$ws.Range("G4").Value = 3.663523847494791
$ws.Range("J4").Value = 10.16385254972159
$ws.Range("K4").Value = 10.72773085843693
$ws.Range("M4").Value = 15.55091928326991
$ws.Range("N4").Value = 19.76701246850477
$ws.Range("O4").Value = 25.67194038360059
$ws.Range("B5").Value = 7.687669738781905
$ws.Range("D5").Value = 8.926913706672766
$ws.Range("E5").Value = 13.61055079596423
$ws.Range("F5").Value = 33.92517398455471
$ws.Range("G5").Value = 3.664061789603944
$ws.Range("J5").Value = 10.17021259721546
$ws.Range("K5").Value = 10.65680099538089
$ws.Range("M5").Value = 15.52732386085421
$ws.Range("N5").Value = 19.78354527948037
$ws.Range("O5").Value = 25.68996536261425
$ws.Range("B6").Value = 7.684913624505908
$ws.Range("D6").Value = 8.926366763637622
$ws.Range("E6").Value = 13.61116720033214
$ws.Range("F6").Value = 33.92772667680002
$ws.Range("G6").Value = 3.664152110002129
$ws.Range("J6").Value = 10.17128590154308
$ws.Range("K6").Value = 10.64500040711654
$ws.Range("M6").Value = 15.5234398943495
$ws.Range("N6").Value = 19.78631850379194
$ws.Range("O6").Value = 25.69303937419146
$ws.Range("B7").Value = 7.704188944280024
$ws.Range("D7").Value = 8.930282644808818
$ws.Range("E7").Value = 13.60708710313583
$ws.Range("F7").Value = 33.91069711253405
$ws.Range("G7").Value = 3.663531035665135
$ws.Range("J7").Value = 10.16393716912049
$ws.Range("K7").Value = 10.72677586595305
$ws.Range("M7").Value = 15.55059879748341
$ws.Range("N7").Value = 19.7672335619871
$ws.Range("O7").Value = 25.67217804452192
$ws.Range("B8").Value = 7.792701729515643
$ws.Range("D8").Value = 8.950500825297345
$ws.Range("E8").Value = 13.59398876235212
$ws.Range("F8").Value = 33.85261383253562
$ws.Range("G8").Value = 3.660934717919423
$ws.Range("J8").Value = 10.13401018914231
$ws.Range("K8").Value = 11.08360350953656
$ws.Range("M8").Value = 15.67525930920266
$ws.Range("N8").Value = 19.68707143736001
$ws.Range("O8").Value = 25.59190943682517
$ws.Range("B9").Value = 7.977404214470302
$ws.Range("D9").Value = 9.000514882165806
$ws.Range("E9").Value = 13.58589057553301
$ws.Range("F9").Value = 33.80015958769565
$ws.Range("G9").Value = 3.656358685032093
$ws.Range("J9").Value = 10.08429330820655
$ws.Range("K9").Value = 11.76238454201373
$ws.Range("M9").Value = 15.93569688156404
$ws.Range("N9").Value = 19.54436070598671
$ws.Range("O9").Value = 25.47702172187866
$ws.Range("B10").Value = 8.118686304804127
$ws.Range("D10").Value = 9.043389254896939
$ws.Range("E10").Value = 13.59061100581231
$ws.Range("F10").Value = 33.79913124120358
$ws.Range("G10").Value = 3.653307554511063
$ws.Range("J10").Value = 10.05321577621639
$ws.Range("K10").Value = 12.2430422203188
$ws.Range("M10").Value = 16.1355376315477
$ws.Range("N10").Value = 19.44825681824219
$ws.Range("O10").Value = 25.41865690651467
$ws.Range("B11").Value = 8.183908010001156
$ws.Range("D11").Value = 9.064179670881634
$ws.Range("E11").Value = 13.59506083133879
$ws.Range("F11").Value = 33.80680253308995
$ws.Range("G11").Value = 3.65198632310555
$ws.Range("J11").Value = 10.04025710354272
$ws.Range("K11").Value = 12.45679566626435
$ws.Range("M11").Value = 16.22801925060591
$ws.Range("N11").Value = 19.40641744548477
$ws.Range("O11").Value = 25.39778121913746
$ws.Range("B12").Value = 8.208720387687421
$ws.Range("D12").Value = 9.072233299686012
$ws.Range("E12").Value = 13.59707538594883
$ws.Range("F12").Value = 33.81087581822662
$ws.Range("G12").Value = 3.651495551837446
$ws.Range("J12").Value = 10.03551916864199
$ws.Range("K12").Value = 12.53695966163967
$ws.Range("M12").Value = 16.26324195332819
$ws.Range("N12").Value = 19.39084278071894
$ws.Range("O12").Value = 25.39069337206876
$ws.Range("B13").Value = 8.203371871338332
$ws.Range("D13").Value = 9.070490837874221
$ws.Range("E13").Value = 13.5966268882107
$ws.Range("F13").Value = 33.80994664356728
$ws.Range("G13").Value = 3.65160082420228
$ws.Range("J13").Value = 10.0365320450949
$ws.Range("K13").Value = 12.5197307227288
$ws.Range("M13").Value = 16.25564752293812
$ws.Range("N13").Value = 19.39418511813229
$ws.Range("O13").Value = 25.39218349870253
$ws.Range("B14").Value = 8.185947176701468
$ws.Range("D14").Value = 9.064838650454893
$ws.Range("E14").Value = 13.59521997602649
$ws.Range("F14").Value = 33.80711424448003
$ws.Range("G14").Value = 3.651945755923343
$ws.Range("J14").Value = 10.03986392057242
$ws.Range("K14").Value = 12.46340684292418
$ws.Range("M14").Value = 16.23091313676675
$ws.Range("N14").Value = 19.40513072396195
$ws.Range("O14").Value = 25.39718171098522
$ws.Range("B15").Value = 8.175288269814573
$ws.Range("D15").Value = 9.061399931206401
$ws.Range("E15").Value = 13.5944010621709
$ws.Range("F15").Value = 33.80553139367826
$ws.Range("G15").Value = 3.652158278910065
$ws.Range("J15").Value = 10.04192682204809
$ws.Range("K15").Value = 12.42880316414212
$ws.Range("M15").Value = 16.21578818194049
$ws.Range("N15").Value = 19.41187022062138
$ws.Range("O15").Value = 25.40034973316503
$ws.Range("B16").Value = 8.114441138489809
$ws.Range("D16").Value = 9.042056082731206
$ws.Range("E16").Value = 13.59036639926829
$ws.Range("F16").Value = 33.79879352767227
$ws.Range("G16").Value = 3.653395239089846
$ws.Range("J16").Value = 10.05408635164828
$ws.Range("K16").Value = 12.22896775814149
$ws.Range("M16").Value = 16.12952332471789
$ws.Range("N16").Value = 19.45102882519765
$ws.Range("O16").Value = 25.42013549051248
$ws.Range("B17").Value = 8.077341131700091
$ws.Range("D17").Value = 9.030515827885104
$ws.Range("E17").Value = 13.58847991170551
$ws.Range("F17").Value = 33.79674349569981
$ws.Range("G17").Value = 3.654171134918582
$ws.Range("J17").Value = 10.06184752913107
$ws.Range("K17").Value = 12.10506469084483
$ws.Range("M17").Value = 16.07698834124104
$ws.Range("N17").Value = 19.4755317236391
$ws.Range("O17").Value = 25.4337278753295
$ws.Range("B18").Value = 8.056093125291955
$ws.Range("D18").Value = 9.023999552784895
$ws.Range("E18").Value = 13.58761168124578
$ws.Range("F18").Value = 33.79633063987345
$ws.Range("G18").Value = 3.654623694657443
$ws.Range("J18").Value = 10.06642250889115
$ws.Range("K18").Value = 12.03334179069622
$ws.Range("M18").Value = 16.04692114613004
$ws.Range("N18").Value = 19.48980205011995
$ws.Range("O18").Value = 25.44207993222697
$ws.Range("B19").Value = 8.048915220626141
$ws.Range("D19").Value = 9.021814226849461
$ws.Range("E19").Value = 13.5873549990255
$ws.Range("F19").Value = 33.79632250969325
$ws.Range("G19").Value = 3.654778004504514
$ws.Range("J19").Value = 10.06799058294425
$ws.Range("K19").Value = 12.00898152922669
$ws.Range("M19").Value = 16.03676732905009
$ws.Range("N19").Value = 19.49466415727205
$ws.Range("O19").Value = 25.44499947879679
$ws.Range("B20").Value = 8.081281242514073
$ws.Range("D20").Value = 9.031731776214432
$ws.Range("E20").Value = 13.58865830165046
$ws.Range("F20").Value = 33.7968824313067
$ws.Range("G20").Value = 3.654087889353843
$ws.Range("J20").Value = 10.06100985688039
$ws.Range("K20").Value = 12.11830225977383
$ws.Range("M20").Value = 16.08256547833913
$ws.Range("N20").Value = 19.47290504750936
$ws.Range("O20").Value = 25.43222565742691
$ws.Range("B21").Value = 8.19106230821993
$ws.Range("D21").Value = 9.066493964450279
$ws.Range("E21").Value = 13.59562429125836
$ws.Range("F21").Value = 33.80791450127462
$ws.Range("G21").Value = 3.651844182295073
$ws.Range("J21").Value = 10.03888067708515
$ws.Range("K21").Value = 12.47997226310109
$ws.Range("M21").Value = 16.23817294027889
$ws.Range("N21").Value = 19.40190844291573
$ws.Range("O21").Value = 25.39569142415172
$ws.Range("B22").Value = 8.26346661549228
$ws.Range("D22").Value = 9.090264489092249
$ws.Range("E22").Value = 13.60209668836443
$ws.Range("F22").Value = 33.82193296161682
$ws.Range("G22").Value = 3.650433435381719
$ws.Range("J22").Value = 10.02540427943225
$ws.Range("K22").Value = 12.71176879791262
$ws.Range("M22").Value = 16.34103693802032
$ws.Range("N22").Value = 19.35707545021506
$ws.Range("O22").Value = 25.37657855594094
$ws.Range("B23").Value = 8.224770470895308
$ws.Range("D23").Value = 9.077482969495344
$ws.Range("E23").Value = 13.59846716804242
$ws.Range("F23").Value = 33.81382898690679
$ws.Range("G23").Value = 3.651181301499526
$ws.Range("J23").Value = 10.03250672487554
$ws.Range("K23").Value = 12.58849652263102
$ws.Range("M23").Value = 16.28603792429823
$ws.Range("N23").Value = 19.38086062761206
$ws.Range("O23").Value = 25.38634316460219
$ws.Range("B24").Value = 8.079499663524036
$ws.Range("D24").Value = 9.031181676951554
$ws.Range("E24").Value = 13.58857697758125
$ws.Range("F24").Value = 33.79681723314228
$ws.Range("G24").Value = 3.654125504464011
$ws.Range("J24").Value = 10.06138821657666
$ws.Range("K24").Value = 12.11231907336534
$ws.Range("M24").Value = 16.08004362959164
$ws.Range("N24").Value = 19.47409199677067
$ws.Range("O24").Value = 25.43290313558618
$ws.Range("B25").Value = 7.926368505068569
$ws.Range("D25").Value = 8.985893741604112
$ws.Range("E25").Value = 13.58620310285966
$ws.Range("F25").Value = 33.80776314270908
$ws.Range("G25").Value = 3.657541793479905
$ws.Range("J25").Value = 10.09678469706365
$ws.Range("K25").Value = 11.58156708600009
$ws.Range("M25").Value = 15.93569688156404
$ws.Range("N25").Value = 19.58142574758202
$ws.Range("O25").Value = 25.50353624725232
